$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 ---
$ws.Range("A15").Value = 130870831
$ws.Range("B15").Value = 83089
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 1312
$ws.Range("F15").Value = "Gammelgransskål"
$ws.Range("G15").Value = "Pseudographis pinicola"
$ws.Range("H15").Value = "(Nyl.) Rehm"
$ws.Range("Q15").Value = 583209
$ws.Range("R15").Value = 6959416

# --- Row 16 ---
$ws.Range("A16").Value = 130870818
$ws.Range("B16").Value = 92267
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 1209
$ws.Range("F16").Value = "Rynkskinn"
$ws.Range("G16").Value = "Hermanssonia centrifuga"
$ws.Range("H16").Value = "(P. Karst.) Zmitr."
$ws.Range("Q16").Value = 583241
$ws.Range("R16").Value = 6959405

# --- Row 17 ---
$ws.Range("A17").Value = 130870792
$ws.Range("B17").Value = 91808
$ws.Range("E17").Value = 1202
$ws.Range("F17").Value = "Ullticka"
$ws.Range("G17").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H17").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q17").Value = 583131
$ws.Range("R17").Value = 6959482

# --- Row 18 ---
$ws.Range("A18").Value = 130870816
$ws.Range("B18").Value = 80377
$ws.Range("E18").Value = 6462
$ws.Range("F18").Value = "Stuplav"
$ws.Range("G18").Value = "Nephroma bellum"
$ws.Range("H18").Value = "(Spreng.) Tuck."
$ws.Range("Q18").Value = 582711
$ws.Range("R18").Value = 6959664

# --- Row 19 ---
$ws.Range("A19").Value = 130870817
$ws.Range("B19").Value = 91819
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 1205
$ws.Range("F19").Value = "Stor aspticka"
$ws.Range("G19").Value = "Phellinus populicola"
$ws.Range("H19").Value = "Niemelä"
$ws.Range("Q19").Value = 582663
$ws.Range("R19").Value = 6959537
$ws.Range("AC19").ClearContents()

# --- Row 20 ---
$ws.Range("A20").Value = 130870795
$ws.Range("B20").Value = 57884
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("M20").Value = "äldre spår"
$ws.Range("Q20").Value = 583098
$ws.Range("R20").Value = 6959481
$ws.Range("AC20").Value = "Äldre ringhack på tall"
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("N20").Value = ""

# --- Row 21 ---
$ws.Range("A21").Value = 130870823
$ws.Range("B21").Value = 79243
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = "Garnlav"
$ws.Range("G21").Value = "Alectoria sarmentosa"
$ws.Range("H21").Value = "(Ach.) Ach."
$ws.Range("Q21").Value = 582529
$ws.Range("R21").Value = 6959663
$ws.Range("AC21").Value = "Med apothecier"
$ws.Range("M21").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("N21").ClearContents()

# --- Row 27 ---
$ws.Range("A27").Value = 130870824
$ws.Range("Q27").Value = 582685
$ws.Range("R27").Value = 6959786

# --- Row 28 ---
$ws.Range("A28").Value = 130870815
$ws.Range("B28").Value = 57884
$ws.Range("E28").Value = 100109
$ws.Range("F28").Value = "Tretåig hackspett"
$ws.Range("G28").Value = "Picoides tridactylus"
$ws.Range("H28").Value = "(Linnaeus, 1758)"
$ws.Range("M28").Value = "färska spår"
$ws.Range("Q28").Value = 583170
$ws.Range("R28").Value = 6959447
$ws.Range("AC28").Value = "Färska och äldre ringhack på tall"
$ws.Range("K28").Value = ""
$ws.Range("L28").Value = ""
$ws.Range("N28").Value = ""

# --- Row 29 ---
$ws.Range("A29").Value = 130870825
$ws.Range("B29").Value = 79243
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("Q29").Value = 582799
$ws.Range("R29").Value = 6959667
$ws.Range("M29").ClearContents()
$ws.Range("AC29").ClearContents()
$ws.Range("K29").ClearContents()
$ws.Range("L29").ClearContents()
$ws.Range("N29").ClearContents()

# --- Row 30 ---
$ws.Range("A30").Value = 130870827
$ws.Range("B30").Value = 79243
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("Q30").Value = 583142
$ws.Range("R30").Value = 6959494

# --- Row 31 ---
$ws.Range("A31").Value = 130870804
$ws.Range("B31").Value = 57884
$ws.Range("E31").Value = 100109
$ws.Range("F31").Value = "Tretåig hackspett"
$ws.Range("G31").Value = "Picoides tridactylus"
$ws.Range("H31").Value = "(Linnaeus, 1758)"
$ws.Range("M31").Value = "äldre spår"
$ws.Range("Q31").Value = 582667
$ws.Range("R31").Value = 6959804
$ws.Range("AC31").Value = "Äldre ringhack på tall"
$ws.Range("K31").Value = ""
$ws.Range("L31").Value = ""
$ws.Range("N31").Value = ""

# --- Row 32 ---
$ws.Range("A32").Value = 130870799
$ws.Range("M32").Value = "färska spår"
$ws.Range("Q32").Value = 582540
$ws.Range("R32").Value = 6959611
$ws.Range("AC32").Value = "Färska och äldre ringhack på tall. Hela stammen full."

# --- Row 33 ---
$ws.Range("A33").Value = 130870794
$ws.Range("B33").Value = 91808
$ws.Range("E33").Value = 1202
$ws.Range("F33").Value = "Ullticka"
$ws.Range("G33").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H33").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q33").Value = 583237
$ws.Range("R33").Value = 6959408
$ws.Range("M33").ClearContents()
$ws.Range("AC33").ClearContents()
$ws.Range("K33").ClearContents()
$ws.Range("L33").ClearContents()
$ws.Range("N33").ClearContents()

# --- Row 37 ---
$ws.Range("A37").Value = 130870810
$ws.Range("B37").Value = 57884
$ws.Range("E37").Value = 100109
$ws.Range("F37").Value = "Tretåig hackspett"
$ws.Range("G37").Value = "Picoides tridactylus"
$ws.Range("H37").Value = "(Linnaeus, 1758)"
$ws.Range("M37").Value = "äldre spår"
$ws.Range("Q37").Value = 582825
$ws.Range("R37").Value = 6959676
$ws.Range("AC37").Value = "Äldre ringhack på tall"
$ws.Range("K37").Value = ""
$ws.Range("L37").Value = ""
$ws.Range("N37").Value = ""

# --- Row 38 ---
$ws.Range("A38").Value = 130870808
$ws.Range("M38").Value = "färska spår"
$ws.Range("Q38").Value = 582781
$ws.Range("R38").Value = 6959717
$ws.Range("AC38").Value = "Färska ringhack på tall"

# --- Row 39 ---
$ws.Range("A39").Value = 130870791
$ws.Range("B39").Value = 91808
$ws.Range("E39").Value = 1202
$ws.Range("F39").Value = "Ullticka"
$ws.Range("G39").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H39").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q39").Value = 582769
$ws.Range("R39").Value = 6959717
$ws.Range("M39").ClearContents()
$ws.Range("AC39").ClearContents()
$ws.Range("K39").ClearContents()
$ws.Range("L39").ClearContents()
$ws.Range("N39").ClearContents()
